$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WBS")

# Delete rows 33-35 (old tasks "Cài đặt module", "kiểm thử module", "Tích hợp hệ thống")
# This shifts rows 36-40 up to become rows 33-37
$ws.Rows("33:35").Delete()

# Update % completed (column F) values for tasks that were completed in week 8 & 9
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = 1

$ws.Range("F20").Value = 1
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("F24").Value = 1
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = 0.7
$ws.Range("F27").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("F31").Value = 0.8
$ws.Range("F32").Value = 1

# Update software version note
$ws.Range("G20").Value = "Version: 4.0"

# Move selection (cosmetic, matches author's cursor position after edit)
$ws.Range("G14").Select()
